# Insert a new weekly price row at row 7 (pushes the existing rows 7:103 down
# to 8:104, matching the diff's dimension change A1:R103 -> A1:R104 and the
# "shift every row down by one" pattern visible for rows 8-104).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new week's data.
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Terminal La Palmera de La Serena"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = "2023-02-06"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 100114007
$ws.Range("G7").Value = "Jengibre"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 23000
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 23500
$ws.Range("N7").Value = "$/caja 13 kilos"
$ws.Range("O7").Value = "Perú"
$ws.Range("P7").Value = 1808
$ws.Range("Q7").Value = 13
$ws.Range("R7").Value = "Hortaliza"
